$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44421
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 15400
$ws.Range("O2").Value = "Provincia de Limarí"
$ws.Range("P2").Value = 513

# Row 3
$ws.Range("D3").Value = 44446
$ws.Range("J3").Value = 25
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 14000
$ws.Range("M3").Value = 14000
$ws.Range("P3").Value = 467

# Row 5
$ws.Range("D5").Value = 44435
$ws.Range("J5").Value = 25
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 14000
$ws.Range("M5").Value = 14000
$ws.Range("P5").Value = 467

# Row 6
$ws.Range("O6").Value = "Provincia del Elquí"

# Row 7
$ws.Range("D7").Value = 44425
$ws.Range("J7").Value = 35
$ws.Range("O7").Value = "Provincia de Limarí"

# Row 8
$ws.Range("D8").Value = 44376
$ws.Range("J8").Value = 25
$ws.Range("K8").Value = 18000
$ws.Range("L8").Value = 18000
$ws.Range("M8").Value = 18000
$ws.Range("P8").Value = 600

# Row 9
$ws.Range("D9").Value = 44418
$ws.Range("J9").Value = 30
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 15000
$ws.Range("P9").Value = 500

# Row 10
$ws.Range("D10").Value = 44467
$ws.Range("J10").Value = 35

# Row 11
$ws.Range("D11").Value = 44432
$ws.Range("O11").Value = "Provincia del Elquí"

# Row 12
$ws.Range("D12").Value = 44449
$ws.Range("J12").Value = 45
$ws.Range("K12").Value = 12000
$ws.Range("L12").Value = 12000
$ws.Range("M12").Value = 12000
$ws.Range("P12").Value = 400

# Row 13
$ws.Range("D13").Value = 44474
$ws.Range("J13").Value = 45
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 10000
$ws.Range("M13").Value = 10000
$ws.Range("P13").Value = 333

# Row 14
$ws.Range("D14").Value = 44460
$ws.Range("J14").Value = 45
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 13000
$ws.Range("M14").Value = 13000
$ws.Range("P14").Value = 433
